$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") and column E ("Volume(1h)") values refreshed by the
# scheduled GitHub Actions crypto-price update job.
# D values are prefixed with a literal apostrophe (PowerShell: '' -> ')
# so Excel stores them as text instead of auto-coercing numeric-looking
# strings (e.g. "1.00", "34.594.36") into numbers/dates.
$ws.Range("D2").Value = '''34.594.36'
$ws.Range("E2").Value = '  +0.31%  '
$ws.Range("D3").Value = '''1.811.86'
$ws.Range("E3").Value = '  +0.47%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").Value = '''225.85'
$ws.Range("E5").Value = '  -1.24%  '
$ws.Range("E6").Value = '  +3.36%  '
$ws.Range("E7").Value = '  -0.19%  '
$ws.Range("D8").Value = '''37.14'
$ws.Range("E8").Value = '  +6.54%  '
$ws.Range("E9").Value = '  -2.42%  '
$ws.Range("D10").Value = '''0.0684'
$ws.Range("E10").Value = '  -1.56%  '
$ws.Range("E11").Value = '  +1.56%  '
$ws.Range("D12").Value = '''2.075.57'
$ws.Range("E12").Value = '  +0.55%  '
$ws.Range("D13").Value = '''11.35'
$ws.Range("E13").Value = '  +1.23%  '
$ws.Range("D14").Value = '''1.817.90'
$ws.Range("E14").Value = '  +0.89%  '
$ws.Range("D15").Value = '''0.634'
$ws.Range("E15").Value = '  -1.46%  '
$ws.Range("D16").Value = '''34.544.11'
$ws.Range("E16").Value = '  +0.25%  '
$ws.Range("D17").Value = '''4.44'
$ws.Range("E17").Value = '  +1.83%  '
$ws.Range("D18").Value = '''68.84'
$ws.Range("E18").Value = '  -0.36%  '
$ws.Range("D19").Value = '''243.41'
$ws.Range("E19").Value = '  -0.78%  '
$ws.Range("D20").Value = '''0.0₃0779'
$ws.Range("E20").Value = '  -2.32%  '
$ws.Range("D21").Value = '''11.26'
$ws.Range("E21").Value = '  -2.19%  '
$ws.Range("E22").Value = '  -0.19%  '
$ws.Range("E23").Value = '  -1.17%  '
$ws.Range("D24").Value = '''2.21'
$ws.Range("E24").Value = '  +4.42%  '
$ws.Range("D25").Value = '''171.94'
$ws.Range("E25").Value = '  -1.08%  '
$ws.Range("D26").Value = '''7.88'
$ws.Range("E26").Value = '  +0.60%  '
$ws.Range("D27").Value = '''17.27'
$ws.Range("E27").Value = '  +2.71%  '
$ws.Range("E28").Value = '  +1.99%  '
$ws.Range("D29").Value = '''1.00'
$ws.Range("E29").Value = '  -0.19%  '
$ws.Range("D30").Value = '''3.83'
$ws.Range("E30").Value = '  -0.16%  '
$ws.Range("D31").Value = '''3.95'
$ws.Range("E31").Value = '  -1.70%  '
$ws.Range("D32").Value = '''1.23'
$ws.Range("E32").Value = '  -1.27%  '
$ws.Range("D33").Value = '''0.0517'
$ws.Range("E33").Value = '  -2.77%  '
$ws.Range("E34").Value = '  -1.01%  '
$ws.Range("D35").Value = '''1.365.32'
$ws.Range("E35").Value = '  -2.12%  '
$ws.Range("E36").Value = '  -4.16%  '
$ws.Range("E37").Value = '  +0.32%  '
$ws.Range("E38").Value = '  -4.79%  '
$ws.Range("D39").Value = '''0.0187'
$ws.Range("E39").Value = '  -1.52%  '
$ws.Range("E40").Value = '  +0.95%  '
$ws.Range("D41").Value = '''81.39'
$ws.Range("E41").Value = '  -2.48%  '
$ws.Range("E42").Value = '  -1.64%  '
$ws.Range("D43").Value = '''0.941'
$ws.Range("E43").Value = '  -0.89%  '
$ws.Range("E44").Value = '  +5.14%  '
$ws.Range("D46").Value = '''0.0503'
$ws.Range("E46").Value = '  -1.52%  '
$ws.Range("D47").Value = '''1.975.00'
$ws.Range("E48").Value = '  -2.33%  '
$ws.Range("E49").Value = '  -0.20%  '
$ws.Range("D50").Value = '''102.78'
$ws.Range("E50").Value = '  -2.03%  '
$ws.Range("D51").Value = '''0.0₆0123'
$ws.Range("E51").Value = '  -5.43%  '
